$wb = $excel.ActiveWorkbook

$jadc = $wb.Worksheets.Item(" JADC (2022)")
$jadc.Name = "JADC (2022)"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)

$gdpUpper = $wb.Worksheets.Add($null, $last)
$gdpUpper.Name = "GDP upper"
$gdpUpper.Range("A1").Value = "year"
$gdpUpper.Range("B1").Value = "GDP "
$gdpUpper.Range("A3").Value = "https://tntcat.iiasa.ac.at/SspDb/dsd?Action=htmlpage&page=30"

$gdpLower = $wb.Worksheets.Add($null, $gdpUpper)
$gdpLower.Name = "GDP lower"
$gdpLower.Range("A1").Value = "year"
$gdpLower.Range("B1").Value = "GDP "

[void]$gdpLower.Range("D5").Select()
[void]$gdpUpper.Activate()
[void]$gdpUpper.Range("E9").Select()
